$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: insert "nodes" in B1, shift the old headers to C1/D1.
$ws.Range("B1").Value = "nodes"
$ws.Range("C1").Value = "penalty of flows"
$ws.Range("D1").Value = "demand of nodes"

# Give the new D1 header the same (bold/bordered/centered) style as the
# other header cells by copying the style from C1.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null

# The A-column index cells (A2:A10) use the same header style too (s="1").
$ws.Range("C1").Copy() | Out-Null
$ws.Range("A2:A10").PasteSpecial(-4122) | Out-Null

$data = @(
    @(0, "(1,2, p = 15, c = 500)", 15, 0),
    @(1, "(1,4, p = 17, c = 200)", 17, 0),
    @(2, "(2,3, p = 8, c = 100)", 8, 0),
    @(3, "(2,5, p = 7, c = 133)", 7, 0),
    @(4, "(3,6, p = 12, c = 100)", 12, 0),
    @(5, "(4,5, p = 15, c = 700)", 15, 0),
    @(6, "(4,7, p = 1, c = 500)", 1, 0),
    @(7, "(5,6, p = 5, c = 200)", 5, 0),
    @(8, "(5,8, p = 7, c = 500)", 7, -400)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
